$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 14499.75
$ws.Range("I64").Value = 16666.334
$ws.Range("J64").Value = 8000
$ws.Range("K64").Value = 16666.334
$ws.Range("L64").Value = 8000
$ws.Range("M64").Value = -16418.334
$ws.Range("N64").Value = -8496
$ws.Range("H67").Value = 14499.75
$ws.Range("I67").Value = 16666.334
$ws.Range("J67").Value = 8000
$ws.Range("K67").Value = 16666.334
$ws.Range("L67").Value = 8000
$ws.Range("M67").Value = -15808.334
$ws.Range("N67").Value = -9716
$ws.Range("H100").Value = 2691.2
$ws.Range("I100").Value = 2528.4
$ws.Range("J100").Value = 3179.6
$ws.Range("K100").Value = 2528.4
$ws.Range("L100").Value = 3179.6
$ws.Range("M100").Value = -1987.4
$ws.Range("N100").Value = -4261.6
$ws.Range("H132").Value = 15874426
$ws.Range("I132").Value = 16668042
$ws.Range("K132").Value = 50004126
$ws.Range("M132").Value = -50001596
$ws.Range("H135").Value = 5556078
$ws.Range("I135").Value = 469.82053
$ws.Range("J135").Value = 41667530
$ws.Range("K135").Value = 4228.384770000001
$ws.Range("L135").Value = 375007770
$ws.Range("M135").Value = -1693.384770000001
$ws.Range("N135").Value = -375012840
$ws.Range("H137").Value = 4082.7646
$ws.Range("I137").Value = 3967.9333
$ws.Range("J137").Value = 4944
$ws.Range("K137").Value = 11903.7999
$ws.Range("L137").Value = 14832
$ws.Range("M137").Value = -9353.7999
$ws.Range("N137").Value = -19932
$ws.Range("H138").Value = 3807.9534
$ws.Range("I138").Value = 997.6667
$ws.Range("J138").Value = 8550.3125
$ws.Range("K138").Value = 2993.0001
$ws.Range("L138").Value = 25650.9375
$ws.Range("M138").Value = 2146.9999
$ws.Range("N138").Value = -35930.9375

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6717.55
$ws.Range("I32").Value = 5720.952
$ws.Range("J32").Value = 9042.944
$ws.Range("K32").Value = 5720.952
$ws.Range("L32").Value = 9042.944
$ws.Range("M32").Value = -5433.952
$ws.Range("N32").Value = -9616.944
$ws.Range("H61").Value = 17547042
$ws.Range("I61").Value = 22225342
$ws.Range("J61").Value = 3415.3333
$ws.Range("K61").Value = 22225342
$ws.Range("L61").Value = 3415.3333
$ws.Range("M61").Value = -22225130
$ws.Range("N61").Value = -3839.3333
$ws.Range("H74").Value = 35755444
$ws.Range("I74").Value = 40045924
$ws.Range("J74").Value = 1466.6666
$ws.Range("K74").Value = 40045924
$ws.Range("L74").Value = 1466.6666
$ws.Range("M74").Value = -40045050
$ws.Range("N74").Value = -3214.6666
$ws.Range("H77").Value = 35755444
$ws.Range("I77").Value = 40045924
$ws.Range("J77").Value = 1466.6666
$ws.Range("K77").Value = 200229620
$ws.Range("L77").Value = 7333.333000000001
$ws.Range("M77").Value = -200225252
$ws.Range("N77").Value = -16069.333
$ws.Range("H132").Value = 27852552
$ws.Range("I132").Value = 9186.241
$ws.Range("J132").Value = 143203650
$ws.Range("K132").Value = 27558.723
$ws.Range("L132").Value = 429610950
$ws.Range("M132").Value = -25028.723
$ws.Range("N132").Value = -429616010
$ws.Range("H136").Value = 17547042
$ws.Range("I136").Value = 22225342
$ws.Range("J136").Value = 3415.3333
$ws.Range("K136").Value = 66676026
$ws.Range("L136").Value = 10245.9999
$ws.Range("M136").Value = -66673476
$ws.Range("N136").Value = -15345.9999

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 99500
$ws.Range("I87").Value = 99500
$ws.Range("K87").Value = 99500
$ws.Range("M87").Value = -98252
$ws.Range("H90").Value = 99500
$ws.Range("I90").Value = 99500
$ws.Range("K90").Value = 298500
$ws.Range("M90").Value = -292260
$ws.Range("H99").Value = 5336.8335
$ws.Range("I99").Value = 1510
$ws.Range("J99").Value = 6102.2
$ws.Range("K99").Value = 1510
$ws.Range("L99").Value = 6102.2
$ws.Range("M99").Value = -12
$ws.Range("N99").Value = -9098.200000000001
$ws.Range("H107").Value = 1755.8334
$ws.Range("I107").Value = 1682.3334
$ws.Range("J107").Value = 1976.3334
$ws.Range("K107").Value = 1682.3334
$ws.Range("L107").Value = 1976.3334
$ws.Range("M107").Value = 237.6666
$ws.Range("N107").Value = -5816.3334
$ws.Range("H134").Value = 5002270
$ws.Range("I134").Value = 5884410
$ws.Range("J134").Value = 3474.6667
$ws.Range("K134").Value = 17653230
$ws.Range("L134").Value = 10424.0001
$ws.Range("M134").Value = -17650695
$ws.Range("N134").Value = -15494.0001

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34487556
$ws.Range("I31").Value = 3089.4285
$ws.Range("J31").Value = 66673056
$ws.Range("K31").Value = 3089.4285
$ws.Range("L31").Value = 66673056
$ws.Range("M31").Value = -2794.4285
$ws.Range("N31").Value = -66673646
$ws.Range("H34").Value = 34487556
$ws.Range("I34").Value = 3089.4285
$ws.Range("J34").Value = 66673056
$ws.Range("K34").Value = 3089.4285
$ws.Range("L34").Value = 66673056
$ws.Range("M34").Value = -2887.4285
$ws.Range("N34").Value = -66673460
$ws.Range("H80").Value = 69999
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 69999
$ws.Range("K80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("M80").Value = 69999
$ws.Range("N80").Value = -72245
$ws.Range("H83").Value = 69999
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 69999
$ws.Range("K83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("M83").Value = 209997
$ws.Range("N83").Value = -221229
$ws.Range("H99").Value = 6889.2144
$ws.Range("I99").Value = 6922.273
$ws.Range("J99").Value = 6768
$ws.Range("K99").Value = 6922.273
$ws.Range("L99").Value = 6768
$ws.Range("M99").Value = -5424.273
$ws.Range("N99").Value = -9764
$ws.Range("H126").Value = 6889.2144
$ws.Range("I126").Value = 6922.273
$ws.Range("J126").Value = 6768
$ws.Range("K126").Value = 20766.819
$ws.Range("L126").Value = 20304
$ws.Range("M126").Value = -18296.819
$ws.Range("N126").Value = -25244
$ws.Range("H132").Value = 62386.8
$ws.Range("I132").Value = 71118.10000000001
$ws.Range("J132").Value = 9999
$ws.Range("K132").Value = 213354.3
$ws.Range("L132").Value = 29997
$ws.Range("M132").Value = -210824.3
$ws.Range("N132").Value = -35057
$ws.Range("H134").Value = 1494.8667
$ws.Range("I134").Value = 1384.0834
$ws.Range("K134").Value = 4152.2502
$ws.Range("M134").Value = -1617.2502

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 215
$ws.Range("I2").Value = 270
$ws.Range("J2").Value = 105
$ws.Range("K2").Value = 1620
$ws.Range("L2").Value = 630
$ws.Range("M2").Value = -1507
$ws.Range("N2").Value = -856

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 11049.583
$ws.Range("I132").Value = 11569
$ws.Range("J132").Value = 10322.4
$ws.Range("K132").Value = 34707
$ws.Range("L132").Value = 30967.2
$ws.Range("M132").Value = -32177
$ws.Range("N132").Value = -36027.2

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 79000
$ws.Range("I92").Value = 79000
$ws.Range("K92").Value = 79000
$ws.Range("M92").Value = -76504
$ws.Range("H132").Value = 62514212
$ws.Range("I132").Value = 10596.546
$ws.Range("J132").Value = 200022160
$ws.Range("K132").Value = 31789.638
$ws.Range("L132").Value = 600066480
$ws.Range("M132").Value = -29259.638
$ws.Range("N132").Value = -600071540
$ws.Range("H136").Value = 1189.38
$ws.Range("I136").Value = 1117.4474
$ws.Range("J136").Value = 1417.1666
$ws.Range("K136").Value = 3352.3422
$ws.Range("L136").Value = 4251.4998
$ws.Range("M136").Value = -802.3422
$ws.Range("N136").Value = -9351.4998

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 186666
$ws.Range("J46").Value = 186666
$ws.Range("L46").Value = 186666
$ws.Range("N46").Value = -187128
$ws.Range("H86").Value = 50015000
$ws.Range("I86").Value = 100000000
$ws.Range("J86").Value = 30000
$ws.Range("K86").Value = 100000000
$ws.Range("L86").Value = 30000
$ws.Range("M86").Value = -99998877
$ws.Range("N86").Value = -32246
$ws.Range("H89").Value = 50015000
$ws.Range("I89").Value = 100000000
$ws.Range("J89").Value = 30000
$ws.Range("K89").Value = 500000000
$ws.Range("L89").Value = 150000
$ws.Range("M89").Value = -499994384
$ws.Range("N89").Value = -161232
$ws.Range("H122").Value = 47668172
$ws.Range("I122").Value = 52685700
$ws.Range("K122").Value = 158057100
$ws.Range("M122").Value = -158054650
$ws.Range("H132").Value = 1680.3846
$ws.Range("I132").Value = 1789.0303
$ws.Range("J132").Value = 1082.8334
$ws.Range("K132").Value = 5367.090899999999
$ws.Range("L132").Value = 3248.5002
$ws.Range("M132").Value = -2837.090899999999
$ws.Range("N132").Value = -8308.5002
$ws.Range("H134").Value = 186666
$ws.Range("J134").Value = 186666
$ws.Range("L134").Value = 559998
$ws.Range("N134").Value = -565068
$ws.Range("H137").Value = 85645
$ws.Range("J137").Value = 85645
$ws.Range("L137").Value = 85645
$ws.Range("N137").Value = -95845
$ws.Range("H139").Value = 199999
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 199999
$ws.Range("K139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("M139").Value = 199999
$ws.Range("N139").Value = -210279
